# Apply the Jan 23 2023 04:31 UTC symbol-list refresh (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these numeric-looking strings as literal
# text (matching the inlineStr cells in the workbook) instead of coercing them to
# real numbers/percentages; resetting the Style afterwards drops the quote-prefix
# flag so the cell format matches the original (unstyled) cells.
$apos = "'"

$ws.Range("D2").Value = $apos + '304.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = $apos + '1.32%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = $apos + '36.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = $apos + '-4.41%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = $apos + '5.031'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = $apos + '1.16%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = $apos + '0.07826'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = $apos + '1.24%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = $apos + '2.168'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = $apos + '-1.03%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = $apos + '7.900'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = $apos + '-1.25%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = $apos + '0.9177'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = $apos + '0.48%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = $apos + '0.09621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = $apos + '3.71%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = $apos + '0.1864'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = $apos + '3.65%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = $apos + '0.08664'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = $apos + '3.07%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = $apos + '0.03491'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = $apos + '-1.50%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = $apos + '0.09918'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = $apos + '-0.22%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = $apos + '0.001436'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = $apos + '-2.89%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = $apos + '0.005677'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = $apos + '-0.87%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = $apos + '3.458'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = $apos + '-0.38%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = $apos + '4.099'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = $apos + '2.57%'
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = $apos + '2.394'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = $apos + '7.69%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = $apos + '0.3424'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = $apos + '-1.11%'
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = $apos + '0.1349'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = $apos + '2.60%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = $apos + '4.774'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = $apos + '5.36%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = $apos + '0.2208'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = $apos + '-0.96%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = $apos + '0.04605'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = $apos + '-1.11%'
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = $apos + '14.95%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = $apos + '0.001232'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = $apos + '0.29%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = $apos + '0.0001402'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = $apos + '7.94%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = $apos + '0.0004755'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = $apos + '0.15%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = $apos + '0.01828'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = $apos + '4.55%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = $apos + '0.04736'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = $apos + '1.07%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = $apos + '0.007389'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = $apos + '-6.08%'
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = $apos + '0.70%'
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = $apos + '0.007752'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = $apos + '1.18%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = $apos + '0.002232'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = $apos + '-2.41%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = $apos + '0.01110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = $apos + '11.55%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = $apos + '0.00006369'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = $apos + '5.16%'
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = $apos + '0.15%'
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = $apos + '0.0005806'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = $apos + '0.10%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = $apos + '23.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = $apos + '169.66%'
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = $apos + '-25.82%'
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = $apos + '0.15%'
$ws.Range("E51").Style = "Normal"
